$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 391). All of them were updated from serial 45182
# (2023-09-13) to serial 45184 (2023-09-15).
$ws.Range("C2:C391").Value = 45184
